$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N..P -> O..Q)
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and update its selection
$ws.Activate() | Out-Null
$ws.Range("T10").Select() | Out-Null
